$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - force text via NumberFormat "@" then ClearFormats
# to avoid Excel auto-converting numeric-looking strings to numbers, while
# keeping cells style-less to match the source workbook formatting.
$priceUpdates = @(
    @{ Addr = "D2"; Value = "36.630.81" },
    @{ Addr = "D3"; Value = "2.048.22" },
    @{ Addr = "D5"; Value = "246.11" },
    @{ Addr = "D7"; Value = "56.90" },
    @{ Addr = "D9"; Value = "63.14" },
    @{ Addr = "D11"; Value = "0.0752" },
    @{ Addr = "D13"; Value = "0.927" },
    @{ Addr = "D14"; Value = "14.45" },
    @{ Addr = "D15"; Value = "2.348.37" },
    @{ Addr = "D16"; Value = "5.44" },
    @{ Addr = "D17"; Value = "2.054.76" },
    @{ Addr = "D18"; Value = "17.87" },
    @{ Addr = "D19"; Value = "36.525.14" },
    @{ Addr = "D20"; Value = "71.99" },
    @{ Addr = "D21"; Value = "0.0₃0859" },
    @{ Addr = "D22"; Value = "237.72" },
    @{ Addr = "D23"; Value = "5.19" },
    @{ Addr = "D27"; Value = "9.31" },
    @{ Addr = "D28"; Value = "164.57" },
    @{ Addr = "D29"; Value = "20.00" },
    @{ Addr = "D31"; Value = "1.22" },
    @{ Addr = "D32"; Value = "5.04" },
    @{ Addr = "D33"; Value = "0.0602" },
    @{ Addr = "D34"; Value = "4.44" },
    @{ Addr = "D35"; Value = "0.0869" },
    @{ Addr = "D38"; Value = "2.20" },
    @{ Addr = "D39"; Value = "5.08" },
    @{ Addr = "D44"; Value = "94.04" },
    @{ Addr = "D45"; Value = "0.0912" },
    @{ Addr = "D46"; Value = "15.99" },
    @{ Addr = "D47"; Value = "1.379.38" },
    @{ Addr = "D48"; Value = "7.45" },
    @{ Addr = "D50"; Value = "2.26" },
    @{ Addr = "D51"; Value = "46.02" }
)

foreach ($u in $priceUpdates) {
    $rng = $ws.Range($u.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}

# Column E (Volume 1h %) updates - these already stay as text because of the
# "%" sign / leading spaces, so a direct .Value assignment is sufficient.
$volumeUpdates = @(
    @{ Addr = "E2"; Value = "  -0.46%  " },
    @{ Addr = "E3"; Value = "  -0.58%  " },
    @{ Addr = "E4"; Value = "  -0.15%  " },
    @{ Addr = "E5"; Value = "  +0.19%  " },
    @{ Addr = "E6"; Value = "  +2.32%  " },
    @{ Addr = "E7"; Value = "  -0.61%  " },
    @{ Addr = "E8"; Value = "  -0.01%  " },
    @{ Addr = "E9"; Value = "  +7.33%  " },
    @{ Addr = "E10"; Value = "  +0.23%  " },
    @{ Addr = "E11"; Value = "  -2.47%  " },
    @{ Addr = "E12"; Value = "  -3.01%  " },
    @{ Addr = "E13"; Value = "  +5.95%  " },
    @{ Addr = "E14"; Value = "  -3.70%  " },
    @{ Addr = "E16"; Value = "  -2.16%  " },
    @{ Addr = "E17"; Value = "  -0.02%  " },
    @{ Addr = "E18"; Value = "  +2.96%  " },
    @{ Addr = "E19"; Value = "  -0.75%  " },
    @{ Addr = "E20"; Value = "  -1.54%  " },
    @{ Addr = "E21"; Value = "  -2.60%  " },
    @{ Addr = "E22"; Value = "  +0.97%  " },
    @{ Addr = "E23"; Value = "  -4.13%  " },
    @{ Addr = "E24"; Value = "  +0.10%  " },
    @{ Addr = "E25"; Value = "  -2.39%  " },
    @{ Addr = "E26"; Value = "  +3.23%  " },
    @{ Addr = "E27"; Value = "  -5.48%  " },
    @{ Addr = "E28"; Value = "  -1.66%  " },
    @{ Addr = "E29"; Value = "  -1.76%  " },
    @{ Addr = "E30"; Value = "  -1.39%  " },
    @{ Addr = "E31"; Value = "  +5.93%  " },
    @{ Addr = "E32"; Value = "  -5.84%  " },
    @{ Addr = "E33"; Value = "  -1.12%  " },
    @{ Addr = "E34"; Value = "  -7.23%  " },
    @{ Addr = "E35"; Value = "  +2.67%  " },
    @{ Addr = "E36"; Value = "  -0.04%  " },
    @{ Addr = "E37"; Value = "  -0.89%  " },
    @{ Addr = "E38"; Value = "  -6.12%  " },
    @{ Addr = "E39"; Value = "  +3.84%  " },
    @{ Addr = "E40"; Value = "  -5.89%  " },
    @{ Addr = "E41"; Value = "  -1.04%  " },
    @{ Addr = "E42"; Value = "  -1.91%  " },
    @{ Addr = "E43"; Value = "  -3.16%  " },
    @{ Addr = "E44"; Value = "  -2.14%  " },
    @{ Addr = "E45"; Value = "  -3.72%  " },
    @{ Addr = "E46"; Value = "  -3.12%  " },
    @{ Addr = "E48"; Value = "  +7.11%  " },
    @{ Addr = "E49"; Value = "  +2.95%  " },
    @{ Addr = "E50"; Value = "  -3.78%  " },
    @{ Addr = "E51"; Value = "  +1.67%  " }
)

foreach ($u in $volumeUpdates) {
    $ws.Range($u.Addr).Value = $u.Value
}
